# Quality pass: remove content hyphens, tweak a few phrases, per commit message.
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# 1. Why 1 answer paragraph
Replace-Text 'pre-acquisition cost structure' 'pre acquisition cost structure'

# 2. Why 2 answer paragraph
Replace-Text 'has a pre-defined Central Finance model' 'has a predefined Central Finance model'
Replace-Text '11 G&A sub-departments' '11 G&A sub departments'

# 3. Evidence paragraph (Why 2)
Replace-Text '$2.37M non-HC OPEX' '$2.37M non HC OPEX'

# 4. Fragmented multi-entity operation
Replace-Text 'fragmented multi-entity operation' 'fragmented multi entity operation'

# 5. Customer impact paragraph
Replace-Text 'is a back-office function' 'is a back office function'

# 6. Target Central Finance Team heading
Replace-Text 'Target Central Finance Team (18 roles, standardized):' 'Target Central Finance Team (18 positions, standardized):'

# 7. Total in-model cost
Replace-Text 'Total in-model cost: $1,200,000/year' 'Total in model cost: $1,200,000/year'

# 8. Target in-model cost
Replace-Text 'Target in-model cost: $1,200,000' 'Target in model cost: $1,200,000'

# 9. Consolidating audit engagements
Replace-Text 'eliminates duplicate entity-level work' 'eliminates duplicate entity level work'

# 10. Unified chart of accounts
Replace-Text 'replaces multiple entity-specific structures' 'replaces multiple entity specific structures'

# 11. F&A is back-office (second occurrence, bullet point)
Replace-Text 'F&A is back-office. Zero' 'F&A is back office. Zero'

# 12. Phase 3 heading
Replace-Text 'Phase 3: Optimization (weeks 13 to 16)' 'Phase 3: Finalize and measure (weeks 13 to 16)'

# 13. Multi-jurisdiction tax bullet
Replace-Text 'Multi-jurisdiction tax: Retain' 'Multi jurisdiction tax: Retain'
Replace-Text 'non-US filings' 'non US filings'

# 14. AI solution: LLM reconciliation agent
Replace-Text 'It auto-matches 100%' 'It auto matches 100%'

# 15. AI solution: Claude Code pipeline
Replace-Text 'flags benchmark non-compliance' 'flags benchmark non compliance'

# 16. Estimated impact paragraph
Replace-Text 'Prevents post-transformation cost drift' 'Prevents post transformation cost drift'

# 17. Week 4 to 8 bullet
Replace-Text 'pilot LLM reconciliation on the highest-volume entity' 'pilot LLM reconciliation on the highest volume entity'

# 18. Claude Code processed paragraph
Replace-Text 'for one sub-department exceeding' 'for one sub department exceeding'
Replace-Text 'aggregating employee and non-employee data' 'aggregating employee and non employee data'

# 19. As VP of Operations paragraph
Replace-Text 'identify the highest-gap function' 'identify the highest gap function'

# 20. The AI tooling makes this repeatable paragraph
Replace-Text 'produces a board-ready deep dive' 'produces a board ready deep dive'
